$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '92.865.77'
$ws.Range('E2').Value = '  +1.14%  '

$ws.Range('D3').Value = '3.116.94'
$ws.Range('E3').Value = '  -1.86%  '

$ws.Range('E4').Value = '  -0.03%  '

$ws.Range('D5').Value = '241.32'
$ws.Range('E5').Value = '  +0.36%  '

$ws.Range('D6').Value = '614.85'
$ws.Range('E6').Value = '  -1.39%  '

$ws.Range('E7').Value = '  -5.14%  '

$ws.Range('E8').Value = '  +7.03%  '

$ws.Range('E9').Value = '  -0.09%  '

$ws.Range('D10').Value = '3.113.59'
$ws.Range('E10').Value = '  -1.86%  '

$ws.Range('D11').Value = '0.729'
$ws.Range('E11').Value = '  -3.08%  '

$ws.Range('E12').Value = '  -2.21%  '

$ws.Range('D13').Value = '0.0000254'
$ws.Range('E13').Value = '  +1.79%  '

$ws.Range('D14').Value = '92.514.22'
$ws.Range('E14').Value = '  +0.91%  '

$ws.Range('D15').Value = '34.40'
$ws.Range('E15').Value = '  -3.46%  '

$ws.Range('D16').Value = '5.49'
$ws.Range('E16').Value = '  -0.79%  '

$ws.Range('D17').Value = '3.695.18'
$ws.Range('E17').Value = '  -1.41%  '

$ws.Range('D18').Value = '3.105.28'
$ws.Range('E18').Value = '  -1.35%  '

$ws.Range('E19').Value = '  -0.71%  '

$ws.Range('D20').Value = '14.72'
$ws.Range('E20').Value = '  -4.95%  '

$ws.Range('E21').Value = '  -2.91%  '

$ws.Range('D22').Value = '9.39'
$ws.Range('E22').Value = '  +1.24%  '

$ws.Range('B23').Value = 'PEPE'
$ws.Range('C23').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D23').Value = '0.0000204'
$ws.Range('E23').Value = '  -3.64%  '

$ws.Range('B24').Value = 'BitcoinCash'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D24').Value = '447.48'
$ws.Range('E24').Value = '  +0.16%  '

$ws.Range('D25').Value = '5.79'
$ws.Range('E25').Value = '  -6.26%  '

$ws.Range('D26').Value = '86.83'
$ws.Range('E26').Value = '  -3.15%  '

$ws.Range('D27').Value = '11.73'
$ws.Range('E27').Value = '  -3.61%  '

$ws.Range('D28').Value = '3.280.19'
$ws.Range('E28').Value = '  -1.31%  '

$ws.Range('E29').Value = '  -0.02%  '

$ws.Range('E30').Value = '  +4.27%  '

$ws.Range('D31').Value = '0.231'
$ws.Range('E31').Value = '  +0.78%  '

$ws.Range('D32').Value = '0.168'
$ws.Range('E32').Value = '  -2.33%  '

$ws.Range('D33').Value = '9.29'
$ws.Range('E33').Value = '  -1.67%  '

$ws.Range('E34').Value = '  +12.80%  '

$ws.Range('D35').Value = '8.07'
$ws.Range('E35').Value = '  +3.85%  '

$ws.Range('E36').Value = '  -1.92%  '

$ws.Range('D37').Value = '4.25'
$ws.Range('E37').Value = '  +9.27%  '

$ws.Range('D38').Value = '26.16'
$ws.Range('E38').Value = '  -1.89%  '

$ws.Range('E39').Value = '  -1.68%  '

$ws.Range('B40').Value = 'Bittensor'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D40').Value = '483.76'
$ws.Range('E40').Value = '  -5.88%  '

$ws.Range('B41').Value = 'Fetch.AI'
$ws.Range('C41').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D41').Value = '1.31'
$ws.Range('E41').Value = '  -3.16%  '

$ws.Range('B42').Value = 'dogwifhat'
$ws.Range('C42').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D42').Value = '3.50'
$ws.Range('E42').Value = '  +0.47%  '

$ws.Range('B43').Value = 'PolygonEcosystemToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D43').Value = '0.438'
$ws.Range('E43').Value = '  -3.48%  '

$ws.Range('D44').Value = '23.10'
$ws.Range('E44').Value = '  +4.27%  '

$ws.Range('E45').Value = '  -0.09%  '

$ws.Range('D46').Value = '161.00'
$ws.Range('E46').Value = '  +2.49%  '

$ws.Range('E47').Value = '  -0.62%  '

$ws.Range('D48').Value = '0.694'
$ws.Range('E48').Value = '  -3.55%  '

$ws.Range('D49').Value = '1.39'
$ws.Range('E49').Value = '  -1.05%  '

$ws.Range('D50').Value = '0.0337'
$ws.Range('E50').Value = '  +0.49%  '

$ws.Range('B51').Value = 'Filecoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D51').Value = '4.41'
$ws.Range('E51').Value = '  -1.53%  '
